# Weekly fruit/vegetable price update: two new price entries were recorded
# for "Ajo" (garlic) at "Vega Modelo de Temuco" and inserted into the table
# right after the existing row for 2022-12-06 (row 1059), pushing every
# subsequent row down by two.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two fresh rows at 1060 (each Insert() pushes the current row 1060
# and everything below it down by one).
$ws.Rows.Item(1060).Insert()
$ws.Rows.Item(1060).Insert()

# New row 1060
$ws.Range("A1060").Value = 10
$ws.Range("B1060").Value = "Vega Modelo de Temuco"
$ws.Range("C1060").Value = "La Araucanía"
$ws.Range("D1060").Value = 45132
$ws.Range("E1060").Value = 9
$ws.Range("F1060").Value = 100112003
$ws.Range("G1060").Value = "Ajo"
$ws.Range("H1060").Value = "Chino"
$ws.Range("I1060").Value = "Primera"
$ws.Range("J1060").Value = 305
$ws.Range("K1060").Value = 20000
$ws.Range("L1060").Value = 23000
$ws.Range("M1060").Value = 20541
$ws.Range("N1060").Value = "$/caja 10 kilos"
$ws.Range("O1060").Value = "China"
$ws.Range("P1060").Value = 2054
$ws.Range("Q1060").Value = 10
$ws.Range("R1060").Value = "Hortaliza"

# New row 1061
$ws.Range("A1061").Value = 10
$ws.Range("B1061").Value = "Vega Modelo de Temuco"
$ws.Range("C1061").Value = "La Araucanía"
$ws.Range("D1061").Value = 45132
$ws.Range("E1061").Value = 9
$ws.Range("F1061").Value = 100112003
$ws.Range("G1061").Value = "Ajo"
$ws.Range("H1061").Value = "Chino"
$ws.Range("I1061").Value = "Primera"
$ws.Range("J1061").Value = 65
$ws.Range("K1061").Value = 26000
$ws.Range("L1061").Value = 26000
$ws.Range("M1061").Value = 26000
$ws.Range("N1061").Value = "$/malla 10 kilos"
$ws.Range("O1061").Value = "China"
$ws.Range("P1061").Value = 2600
$ws.Range("Q1061").Value = 10
$ws.Range("R1061").Value = "Hortaliza"
